# LOQ4239.xlsx rebuild — mirrors the upstream "Build site" regeneration that
# collapsed rows 22:23 into the table and shuffled several labels/values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the two trailing rows (old "Norma de recuperação:" dup + "Bibliografia:"
#    rows 22 and 23) — the rebuilt sheet ends at row 21.
$ws.Rows("22:23").Delete()

# 2) Objetivos: row now shows the first docente responsável instead of the
#    old course-objectives paragraph.
$ws.Range("B10").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Value = "11079086 - Herlandí de Souza Andrade"

# 3) Rows 13-21 get new labels/values (section headers shifted, syllabus /
#    bibliography bodies dropped, docente names and dates redistributed).
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2021"
$ws.Range("C13").Value = "01/01/2021"

$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = $null
$ws.Range("C14").Value = $null

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C15").Value = "11079086 - Herlandí de Souza Andrade"

$ws.Range("A16").Value = "Syllabus:"

$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Value = $null
$ws.Range("C17").Value = $null

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C18").Value = "5840560 - Marco Antonio Carvalho Pereira"

$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Range("C19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."

$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."
$ws.Range("C20").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."

$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."

# 4) Row heights: rows 13/14 gain a 60pt band, 15/16 grow to 120pt, 17 drops
#    back to the default (autofit clears the custom height), 18 shrinks from
#    120 to 60, 19 gains a 60pt band, and 21 grows to 120pt.
$ws.Rows(13).RowHeight = 60
$ws.Rows(14).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(16).RowHeight = 120
$ws.Rows(17).EntireRow.AutoFit()
$ws.Rows(18).RowHeight = 60
$ws.Rows(19).RowHeight = 60
$ws.Rows(21).RowHeight = 120
